$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells retain their original text (inline string) representation
# by forcing Text number format before assigning values, then updating the
# symbol/link/price/volume columns to match the refreshed crypto data feed.

$cells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "B7", "C7", "D7", "E7", "B8", "C8", "D8", "E8", "B9", "C9", "D9", "E9", "B10", "C10", "D10", "E10", "B11", "C11", "D11", "E11", "B12", "C12", "D12", "E12", "B13", "C13", "D13", "E13", "B14", "C14", "D14", "E14", "B15", "C15", "D15", "E15", "B16", "C16", "D16", "E16", "B17", "C17", "D17", "E17", "D18", "E18", "D19", "E19", "E20", "D21", "E21", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D27", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "E47", "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "300.60"
$ws.Range("E2").Value = "-4.37%"
$ws.Range("D3").Value = "35.39"
$ws.Range("E3").Value = "-0.85%"
$ws.Range("D4").Value = "5.051"
$ws.Range("E4").Value = "-0.79%"
$ws.Range("D5").Value = "0.07958"
$ws.Range("E5").Value = "-2.18%"
$ws.Range("D6").Value = "1.902"
$ws.Range("E6").Value = "-9.58%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "4.060"
$ws.Range("E7").Value = "-2.01%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "7.794"
$ws.Range("E8").Value = "-1.91%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9218"
$ws.Range("E9").Value = "-1.14%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1426"
$ws.Range("E10").Value = "38.41%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1896"
$ws.Range("E11").Value = "-1.50%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09156"
$ws.Range("E12").Value = "1.60%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03403"
$ws.Range("E13").Value = "-5.66%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09866"
$ws.Range("E14").Value = "-0.24%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001400"
$ws.Range("E15").Value = "-3.04%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005821"
$ws.Range("E16").Value = "0.28%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.519"
$ws.Range("E17").Value = "1.54%"
$ws.Range("D18").Value = "2.942"
$ws.Range("E18").Value = "1.95%"
$ws.Range("D19").Value = "0.3401"
$ws.Range("E19").Value = "-0.21%"
$ws.Range("E20").Value = "-2.05%"
$ws.Range("D21").Value = "5.072"
$ws.Range("E21").Value = "-0.67%"
$ws.Range("D23").Value = "0.04491"
$ws.Range("E23").Value = "-1.41%"
$ws.Range("D24").Value = "0.001218"
$ws.Range("E24").Value = "-2.43%"
$ws.Range("D25").Value = "0.004773"
$ws.Range("E25").Value = "-0.45%"
$ws.Range("D26").Value = "0.0001234"
$ws.Range("E26").Value = "-1.54%"
$ws.Range("D27").Value = "0.0003007"
$ws.Range("E27").Value = "-33.35%"
$ws.Range("D39").Value = "0.01905"
$ws.Range("E39").Value = "-2.79%"
$ws.Range("D40").Value = "0.04717"
$ws.Range("E40").Value = "-3.36%"
$ws.Range("D41").Value = "0.007387"
$ws.Range("E41").Value = "-3.14%"
$ws.Range("D42").Value = "0.009700"
$ws.Range("E42").Value = "25.68%"
$ws.Range("D43").Value = "0.1323"
$ws.Range("E43").Value = "-4.28%"
$ws.Range("D44").Value = "0.002117"
$ws.Range("E44").Value = "0.42%"
$ws.Range("D45").Value = "0.01038"
$ws.Range("E45").Value = "-11.63%"
$ws.Range("D46").Value = "0.00006257"
$ws.Range("E46").Value = "-6.86%"
$ws.Range("E47").Value = "0.07%"
$ws.Range("E48").Value = "-58.88%"
$ws.Range("D49").Value = "0.001663"
$ws.Range("E49").Value = "-2.45%"
$ws.Range("D50").Value = "0.00002107"
$ws.Range("E50").Value = "0.07%"
$ws.Range("D51").Value = "0.0002007"
$ws.Range("E51").Value = "0.07%"
